$wb = $excel.ActiveWorkbook

$ws = $wb.Worksheets.Item("ALC")
$ws.Range("H33").Value = 119.82353
$ws.Range("I33").Value = 120.9375
$ws.Range("J33").Value = 102
$ws.Range("K33").Value = 120.9375
$ws.Range("L33").Value = 102
$ws.Range("M33").Value = 108.0625
$ws.Range("N33").Value = -560

$ws.Range("H40").Value = 2429.1428
$ws.Range("I40").Value = 2150
$ws.Range("J40").Value = 2540.8
$ws.Range("K40").Value = 2150
$ws.Range("L40").Value = 2540.8
$ws.Range("M40").Value = -1975
$ws.Range("N40").Value = -2890.8

$ws.Range("H43").Value = 8355.444
$ws.Range("I43").Value = 8000
$ws.Range("J43").Value = 8533.166999999999
$ws.Range("K43").Value = 8000
$ws.Range("L43").Value = 8533.166999999999
$ws.Range("M43").Value = -7931
$ws.Range("N43").Value = -8671.166999999999

$ws.Range("H131").Value = 3732.8157
$ws.Range("I131").Value = 303.35
$ws.Range("J131").Value = 7543.3335
$ws.Range("K131").Value = 910.0500000000001
$ws.Range("L131").Value = 22630.0005
$ws.Range("M131").Value = 4129.95
$ws.Range("N131").Value = -32710.0005

$ws.Range("H135").Value = 2626.2632
$ws.Range("I135").Value = 910.46155
$ws.Range("J135").Value = 6343.8335
$ws.Range("K135").Value = 8194.15395
$ws.Range("L135").Value = 57094.5015
$ws.Range("M135").Value = -5659.15395
$ws.Range("N135").Value = -62164.5015

$ws = $wb.Worksheets.Item("ARM")
$ws.Range("H45").Value = 967.4286
$ws.Range("I45").Value = 1093
$ws.Range("J45").Value = 800
$ws.Range("K45").Value = 1093
$ws.Range("L45").Value = 800
$ws.Range("M45").Value = -716
$ws.Range("N45").Value = -1554

$ws.Range("H68").Value = 30099
$ws.Range("J68").Value = 30099
$ws.Range("L68").Value = 30099
$ws.Range("N68").Value = -31721

$ws.Range("H71").Value = 30099
$ws.Range("J71").Value = 30099
$ws.Range("L71").Value = 90297
$ws.Range("N71").Value = -98409

$ws = $wb.Worksheets.Item("BSM")
$ws.Range("H86").Value = 7693781.5
$ws.Range("I86").Value = 11112250
$ws.Range("J86").Value = 2226.75
$ws.Range("K86").Value = 11112250
$ws.Range("L86").Value = 2226.75
$ws.Range("M86").Value = -11111127
$ws.Range("N86").Value = -4472.75

$ws.Range("H89").Value = 7693781.5
$ws.Range("I89").Value = 11112250
$ws.Range("J89").Value = 2226.75
$ws.Range("K89").Value = 55561250
$ws.Range("L89").Value = 11133.75
$ws.Range("M89").Value = -55555634
$ws.Range("N89").Value = -22365.75

$ws.Range("H139").Value = 62666.668
$ws.Range("J139").Value = 62666.668
$ws.Range("L139").Value = 62666.668
$ws.Range("N139").Value = -72946.66800000001

$ws = $wb.Worksheets.Item("CRP")
$ws.Range("H58").Value = 7199345
$ws.Range("I58").Value = 15985096
$ws.Range("J58").Value = 11002.546
$ws.Range("K58").Value = 15985096
$ws.Range("L58").Value = 11002.546
$ws.Range("M58").Value = -15984893
$ws.Range("N58").Value = -11408.546

$ws.Range("H99").Value = 142858200
$ws.Range("I99").Value = 250000900
$ws.Range("K99").Value = 250000900
$ws.Range("M99").Value = -249999402

$ws.Range("H126").Value = 142858200
$ws.Range("I126").Value = 250000900
$ws.Range("K126").Value = 750002700
$ws.Range("M126").Value = -750000230

$ws.Range("H132").Value = 10757899
$ws.Range("I132").Value = 18519568
$ws.Range("J132").Value = 10971.23
$ws.Range("K132").Value = 55558704
$ws.Range("L132").Value = 32913.69
$ws.Range("M132").Value = -55556174
$ws.Range("N132").Value = -37973.69

$ws.Range("H134").Value = 26043418
$ws.Range("I134").Value = 25001900
$ws.Range("J134").Value = 31251000
$ws.Range("K134").Value = 75005700
$ws.Range("L134").Value = 93753000
$ws.Range("M134").Value = -75003165
$ws.Range("N134").Value = -93758070

$ws.Range("H135").Value = 30000
$ws.Range("J135").Value = 30000
$ws.Range("L135").Value = 30000
$ws.Range("N135").Value = -40140

$ws.Range("H136").Value = 7199345
$ws.Range("I136").Value = 15985096
$ws.Range("J136").Value = 11002.546
$ws.Range("K136").Value = 47955288
$ws.Range("L136").Value = 33007.638
$ws.Range("M136").Value = -47952738
$ws.Range("N136").Value = -38107.638

$ws = $wb.Worksheets.Item("CUL")
$ws.Range("H75").Value = 2240
$ws.Range("J75").Value = 2366.6667
$ws.Range("L75").Value = 7100.000100000001
$ws.Range("N75").Value = -9096.000100000001

$ws.Range("H78").Value = 2240
$ws.Range("J78").Value = 2366.6667
$ws.Range("L78").Value = 21300.0003
$ws.Range("N78").Value = -31284.0003

$ws = $wb.Worksheets.Item("GSM")
$ws.Range("H126").Value = 1814.1333
$ws.Range("I126").Value = 1301.7142
$ws.Range("J126").Value = 2262.5
$ws.Range("K126").Value = 3905.1426
$ws.Range("L126").Value = 6787.5
$ws.Range("M126").Value = -1435.1426
$ws.Range("N126").Value = -11727.5

$ws.Range("H132").Value = 23282696
$ws.Range("I132").Value = 38501864
$ws.Range("J132").Value = 6321.353
$ws.Range("K132").Value = 115505592
$ws.Range("L132").Value = 18964.059
$ws.Range("M132").Value = -115503062
$ws.Range("N132").Value = -24024.059

$ws.Range("H141").Value = 29714.334
$ws.Range("J141").Value = 29714.334
$ws.Range("L141").Value = 29714.334
$ws.Range("N141").Value = -40074.334

$ws = $wb.Worksheets.Item("LTW")
$ws.Range("H132").Value = 4880524
$ws.Range("I132").Value = 7145205
$ws.Range("J132").Value = 2749.7693
$ws.Range("K132").Value = 21435615
$ws.Range("L132").Value = 8249.3079
$ws.Range("M132").Value = -21433085
$ws.Range("N132").Value = -13309.3079

$ws.Range("H135").Value = 32366.5
$ws.Range("J135").Value = 32366.5
$ws.Range("L135").Value = 32366.5
$ws.Range("N135").Value = -42506.5

$ws.Range("H136").Value = 2712.413
$ws.Range("I136").Value = 4080.524
$ws.Range("J136").Value = 1563.2
$ws.Range("K136").Value = 12241.572
$ws.Range("L136").Value = 4689.6
$ws.Range("M136").Value = -9691.572
$ws.Range("N136").Value = -9789.6

$ws = $wb.Worksheets.Item("WVR")
$ws.Range("H44").Value = 8059.5
$ws.Range("J44").Value = 8059.5
$ws.Range("L44").Value = 8059.5
$ws.Range("N44").Value = -9167.5

$ws.Range("H113").Value = 482.03125
$ws.Range("I113").Value = 301.08694
$ws.Range("J113").Value = 944.44446
$ws.Range("K113").Value = 903.2608200000001
$ws.Range("L113").Value = 2833.33338
$ws.Range("M113").Value = 1266.73918
$ws.Range("N113").Value = -7173.33338

$ws.Range("H126").Value = 46296948
$ws.Range("I126").Value = 13889408
$ws.Range("J126").Value = 111112024
$ws.Range("K126").Value = 41668224
$ws.Range("L126").Value = 333336072
$ws.Range("M126").Value = -41665754
$ws.Range("N126").Value = -333341012
